# Armenia Premier League 2023-2024 — script update 28-11-2023 20:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rows 60/61 (F:V, i.e. match data excluding Indice/pais/torneio/temporada/data_partida)
#        were re-sorted; swap the two rows' F:V content back.
$r60 = $ws.Range("F60:V60")
$r61 = $ws.Range("F61:V61")
$v60 = $r60.Value2
$v61 = $r61.Value2
$r60.Value2 = $v61
$r61.Value2 = $v60

# --- 2) Append new match row 81 (Shirak Gyumri 1-1 Pyunik Yerevan), reusing the
#        number-format/border styling already used by row 80's cells.
$ws.Range("A80").Copy()
$ws.Range("A81").PasteSpecial(-4122)
$ws.Range("E80").Copy()
$ws.Range("E81").PasteSpecial(-4122)

$ws.Range("A81").Value = 80
$ws.Range("B81").Value = "armenia"
$ws.Range("C81").Value = "premier-league"
$ws.Range("D81").Value = "2023-2024"
$ws.Range("E81").Value = 45258.625
$ws.Range("F81").Value = "Shirak Gyumri"
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = "Pyunik Yerevan"
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 8.88
$ws.Range("K81").Value = "27/11/2023 03:12"
$ws.Range("L81").Value = 11.65
$ws.Range("M81").Value = "28/11/2023 13:46"
$ws.Range("N81").Value = 5.71
$ws.Range("O81").Value = "27/11/2023 03:12"
$ws.Range("P81").Value = 5.9
$ws.Range("Q81").Value = "28/11/2023 13:46"
$ws.Range("R81").Value = 1.23
$ws.Range("S81").Value = "27/11/2023 03:12"
$ws.Range("T81").Value = 1.24
$ws.Range("U81").Value = "28/11/2023 13:26"
$ws.Range("V81").Value = "https://www.betexplorer.com/football/armenia/premier-league/shirak-gyumri-pyunik-yerevan/KrGzgW3b/"
